# Update the cryptocurrency price/volume table to the latest scraped values.
# Each text cell is written with an explicit text NumberFormat so numeric-looking
# strings (e.g. "603.85", "0.0000125") are stored verbatim as text instead of
# being coerced into numbers; ClearFormats() afterwards drops the temporary format
# override so the cell's style index is left exactly as it started (no style/format change).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '67.246.25'
$c.ClearFormats()
$c = $ws.Range("E2")
$c.NumberFormat = "@"
$c.Value = '  -0.42%  '
$c.ClearFormats()

# Row 3
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '3.478.98'
$c.ClearFormats()
$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = '  -1.31%  '
$c.ClearFormats()

# Row 4
$c = $ws.Range("E4")
$c.NumberFormat = "@"
$c.Value = '  +0.13%  '
$c.ClearFormats()

# Row 5
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '603.85'
$c.ClearFormats()
$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value = '  -1.74%  '
$c.ClearFormats()

# Row 6
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '150.59'
$c.ClearFormats()
$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = '  -0.81%  '
$c.ClearFormats()

# Row 7
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '3.478.50'
$c.ClearFormats()
$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value = '  -1.29%  '
$c.ClearFormats()

# Row 8
$c = $ws.Range("E8")
$c.NumberFormat = "@"
$c.Value = '  -0.15%  '
$c.ClearFormats()

# Row 9
$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = '  +0.60%  '
$c.ClearFormats()

# Row 10
$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = '  +2.33%  '
$c.ClearFormats()

# Row 11
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '7.58'
$c.ClearFormats()
$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = '  +6.57%  '
$c.ClearFormats()

# Row 12
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.428'
$c.ClearFormats()
$c = $ws.Range("E12")
$c.NumberFormat = "@"
$c.Value = '  +0.68%  '
$c.ClearFormats()

# Row 13
$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value = '  -2.37%  '
$c.ClearFormats()

# Row 14
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '31.98'
$c.ClearFormats()

# Row 15
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '4.066.83'
$c.ClearFormats()
$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = '  -1.29%  '
$c.ClearFormats()

# Row 16
$c = $ws.Range("B16")
$c.NumberFormat = "@"
$c.Value = 'WrappedBTC'
$c.ClearFormats()
$c = $ws.Range("C16")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$c.ClearFormats()
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '67.183.09'
$c.ClearFormats()
$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value = '  -0.43%  '
$c.ClearFormats()

# Row 17
$c = $ws.Range("B17")
$c.NumberFormat = "@"
$c.Value = 'WrappedEther'
$c.ClearFormats()
$c = $ws.Range("C17")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$c.ClearFormats()
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '3.463.27'
$c.ClearFormats()
$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = '  -1.70%  '
$c.ClearFormats()

# Row 18
$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value = '  -0.89%  '
$c.ClearFormats()

# Row 19
$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = '  +1.18%  '
$c.ClearFormats()

# Row 20
$c = $ws.Range("E20")
$c.NumberFormat = "@"
$c.Value = '  -0.44%  '
$c.ClearFormats()

# Row 21
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '9.83'
$c.ClearFormats()
$c = $ws.Range("E21")
$c.NumberFormat = "@"
$c.Value = '  +3.56%  '
$c.ClearFormats()

# Row 22
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '444.64'
$c.ClearFormats()
$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = '  +0.01%  '
$c.ClearFormats()

# Row 24
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '78.03'
$c.ClearFormats()
$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = '  +0.79%  '
$c.ClearFormats()

# Row 25
$c = $ws.Range("E25")
$c.NumberFormat = "@"
$c.Value = '  -0.11%  '
$c.ClearFormats()

# Row 26
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '3.620.91'
$c.ClearFormats()
$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value = '  -1.20%  '
$c.ClearFormats()

# Row 27
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '0.0000125'
$c.ClearFormats()
$c = $ws.Range("E27")
$c.NumberFormat = "@"
$c.Value = '  -4.69%  '
$c.ClearFormats()

# Row 28
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '8.66'
$c.ClearFormats()
$c = $ws.Range("E28")
$c.NumberFormat = "@"
$c.Value = '  +1.65%  '
$c.ClearFormats()

# Row 29
$c = $ws.Range("E29")
$c.NumberFormat = "@"
$c.Value = '  -3.40%  '
$c.ClearFormats()

# Row 30
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '2.49'
$c.ClearFormats()
$c = $ws.Range("E30")
$c.NumberFormat = "@"
$c.Value = '  -1.11%  '
$c.ClearFormats()

# Row 31
$c = $ws.Range("E31")
$c.NumberFormat = "@"
$c.Value = '  +3.06%  '
$c.ClearFormats()

# Row 32
$c = $ws.Range("E32")
$c.NumberFormat = "@"
$c.Value = '  +3.74%  '
$c.ClearFormats()

# Row 33
$c = $ws.Range("E33")
$c.NumberFormat = "@"
$c.Value = '  +0.20%  '
$c.ClearFormats()

# Row 34
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '25.46'
$c.ClearFormats()
$c = $ws.Range("E34")
$c.NumberFormat = "@"
$c.Value = '  -1.58%  '
$c.ClearFormats()

# Row 35
$c = $ws.Range("E35")
$c.NumberFormat = "@"
$c.Value = '  -1.05%  '
$c.ClearFormats()

# Row 36
$c = $ws.Range("E36")
$c.NumberFormat = "@"
$c.Value = '  -0.14%  '
$c.ClearFormats()

# Row 37
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '3.474.39'
$c.ClearFormats()
$c = $ws.Range("E37")
$c.NumberFormat = "@"
$c.Value = '  -1.17%  '
$c.ClearFormats()

# Row 38
$c = $ws.Range("E38")
$c.NumberFormat = "@"
$c.Value = '  -0.91%  '
$c.ClearFormats()

# Row 39
$c = $ws.Range("E39")
$c.NumberFormat = "@"
$c.Value = '  -0.03%  '
$c.ClearFormats()

# Row 40
$c = $ws.Range("E40")
$c.NumberFormat = "@"
$c.Value = '  +6.50%  '
$c.ClearFormats()

# Row 41
$c = $ws.Range("E41")
$c.NumberFormat = "@"
$c.Value = '  +0.14%  '
$c.ClearFormats()

# Row 42
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '176.78'
$c.ClearFormats()
$c = $ws.Range("E42")
$c.NumberFormat = "@"
$c.Value = '  -0.29%  '
$c.ClearFormats()

# Row 43
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '0.0892'
$c.ClearFormats()
$c = $ws.Range("E43")
$c.NumberFormat = "@"
$c.Value = '  +0.49%  '
$c.ClearFormats()

# Row 44
$c = $ws.Range("E44")
$c.NumberFormat = "@"
$c.Value = '  -0.28%  '
$c.ClearFormats()

# Row 45
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '0.889'
$c.ClearFormats()
$c = $ws.Range("E45")
$c.NumberFormat = "@"
$c.Value = '  +0.69%  '
$c.ClearFormats()

# Row 46
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '29.95'
$c.ClearFormats()
$c = $ws.Range("E46")
$c.NumberFormat = "@"
$c.Value = '  +5.19%  '
$c.ClearFormats()

# Row 47
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '46.42'
$c.ClearFormats()
$c = $ws.Range("E47")
$c.NumberFormat = "@"
$c.Value = '  +2.84%  '
$c.ClearFormats()

# Row 48
$c = $ws.Range("E48")
$c.NumberFormat = "@"
$c.Value = '  +2.86%  '
$c.ClearFormats()

# Row 49
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '2.52'
$c.ClearFormats()
$c = $ws.Range("E49")
$c.NumberFormat = "@"
$c.Value = '  -4.50%  '
$c.ClearFormats()

# Row 50
$c = $ws.Range("E50")
$c.NumberFormat = "@"
$c.Value = '  -0.41%  '
$c.ClearFormats()

# Row 51
$c = $ws.Range("B51")
$c.NumberFormat = "@"
$c.Value = 'TheGraph'
$c.ClearFormats()
$c = $ws.Range("C51")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$c.ClearFormats()
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '0.251'
$c.ClearFormats()
$c = $ws.Range("E51")
$c.NumberFormat = "@"
$c.Value = '  +0.07%  '
$c.ClearFormats()
